$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data table: row number -> values for columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T
$colNames = @('E','G','H','I','J','K','M','N','O','P','Q','R','S','T')
$data = @{
    2 = @(3, 17.85965566666667, 53.578967, 0.4399433554208643, 0.4399433554208643, 3, 3.535386, 10.606158, 0.01988747852527457, 0.01988747852527457, 63.140776608754, 568.2669894787861, 0.008749364033269678, 0.008749364033269676)
    3 = @(3, 17.85965566666667, 53.578967, 0.4399433554208643, 0.4399433554208643, 3, 173.8189136666666, 521.456741, 0.9777772251268709, 0.9777772251268707, 3104.345946440727, 27939.11351796654, 0.4301665932764174, 0.4301665932764174)
    4 = @(3, 17.85965566666667, 53.578967, 0.4399433554208643, 0.4399433554208643, 3, 0.06908833333333333, 0.207265, 0.000388640093475982, 0.0003886400934759819, 1.233893843917222, 11.105044595255, 0.0001709796267749019, 0.0001709796267749018)
    5 = @(3, 17.85965566666667, 53.578967, 0.4399433554208643, 0.4399433554208643, 3, 0.346056, 1.038168, 0.001946656254378565, 0.001946656254378564, 6.180441001383999, 55.62396901245599, 0.0008564184844023174, 0.0008564184844023173)
    6 = @(3, 2.191749666666666, 6.575248999999999, 0.05399016199374809, 0.05399016199374809, 3, 3.535386, 10.606158, 0.01988747852527457, 0.01988747852527457, 7.748681087038, 69.738129783342, 0.00107372818722676, 0.00107372818722676)
    7 = @(3, 2.191749666666666, 6.575248999999999, 0.05399016199374809, 0.05399016199374809, 3, 173.8189136666666, 521.456741, 0.9777772251268709, 0.9777772251268707, 380.9675460892787, 3428.707914803509, 0.05279035077839725, 0.05279035077839724)
    8 = @(3, 2.191749666666666, 6.575248999999999, 0.05399016199374809, 0.05399016199374809, 3, 0.06908833333333333, 0.207265, 0.000388640093475982, 0.0003886400934759819, 0.1514243315538889, 1.362818983985, 0.00002098274160403367, 0.00002098274160403366)
    9 = @(3, 2.191749666666666, 6.575248999999999, 0.05399016199374809, 0.05399016199374809, 3, 0.346056, 1.038168, 0.001946656254378565, 0.001946656254378564, 0.7584681226479999, 6.826213103831999, 0.0001051002865200416, 0.0001051002865200416)
    10 = @(3, 7.079289333333333, 21.237868, 0.1743866937543869, 0.1743866937543869, 3, 3.535386, 10.606158, 0.01988747852527457, 0.01988747852527457, 25.028020399016, 225.252183591144, 0.003468111627134003, 0.003468111627134002)
    11 = @(3, 7.079289333333333, 21.237868, 0.1743866937543869, 0.1743866937543869, 3, 173.8189136666666, 521.456741, 0.9777772251268709, 0.9777772251268707, 1230.514381452021, 11074.62943306819, 0.1705113375182138, 0.1705113375182138)
    12 = @(3, 7.079289333333333, 21.237868, 0.1743866937543869, 0.1743866937543869, 3, 0.06908833333333333, 0.207265, 0.000388640093475982, 0.0003886400934759819, 0.4890963012244444, 4.40186671102, 0.00006777366096167237, 0.00006777366096167235)
    13 = @(3, 7.079289333333333, 21.237868, 0.1743866937543869, 0.1743866937543869, 3, 0.346056, 1.038168, 0.001946656254378565, 0.001946656254378564, 2.449830549536, 22.048474945824, 0.0003394709480773767, 0.0003394709480773766)
    14 = @(3, 13.464658, 40.393974, 0.3316797888310007, 0.3316797888310007, 3, 3.535386, 10.606158, 0.01988747852527457, 0.01988747852527457, 47.60276338798801, 428.424870491892, 0.006596274677644132, 0.006596274677644129)
    15 = @(3, 13.464658, 40.393974, 0.3316797888310007, 0.3316797888310007, 3, 173.8189136666666, 521.456741, 0.9777772251268709, 0.9777772251268707, 2340.412226453192, 21063.71003807873, 0.3243089435538424, 0.3243089435538423)
    16 = @(3, 13.464658, 40.393974, 0.3316797888310007, 0.3316797888310007, 3, 0.06908833333333333, 0.207265, 0.000388640093475982, 0.0003886400934759819, 0.9302507801233334, 8.37225702111, 0.0001289040641353741, 0.000128904064135374)
    17 = @(3, 13.464658, 40.393974, 0.3316797888310007, 0.3316797888310007, 3, 0.346056, 1.038168, 0.001946656254378565, 0.001946656254378564, 4.659525688847999, 41.935731199632, 0.0006456665353788291, 0.0006456665353788291)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $colNames.Length; $i++) {
        $cellRef = "$($colNames[$i])$r"
        $ws.Range($cellRef).Value = $vals[$i]
    }
}
